$d = $word.ActiveDocument

# Locate the paragraph that contains "Registrar profesores ... dictan."
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Registrar profesores*dictan*") {
        $pRange = $p.Range

        # Highlight the whole paragraph (including the paragraph mark) in yellow,
        # matching the surrounding list items' formatting.
        $pRange.Font.HighlightColorIndex = 7  # wdYellow

        # Merge the two separate runs ("Registrar profesores" + " y consultar
        # qué cursos dictan.") into a single run by replacing the paragraph's
        # text in-place.
        $pRange.Find.Execute(
            "Registrar profesores y consultar qué cursos dictan.",
            $false, $false, $false, $false, $false, $true, 1, $false,
            "Registrar profesores y consultar qué cursos dictan.", 2
        ) | Out-Null

        break
    }
}
